# Update fleet operator assignments ("Stato Attuale") and log the two
# operator changes in the change-history sheet ("Storico Passaggi").

$wb = $excel.ActiveWorkbook
$wsStato = $wb.Worksheets.Item("Stato Attuale")
$wsStorico = $wb.Worksheets.Item("Storico Passaggi")

# Helper: write a value into a cell as literal text (so date-looking
# strings like "2026-02-13" are NOT auto-converted to Excel date serials),
# while leaving the cell's style/number-format as plain "Normal" afterwards.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- "Stato Attuale" sheet: locate rows by targa (plate) in column A ---
$usedRange = $wsStato.UsedRange
$lastRow = $usedRange.Rows.Count

$rowGL592TH = 0
$rowGY983FY = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $targa = $wsStato.Cells.Item($r, 1).Value2
    if ($targa -eq "GL592TH") { $rowGL592TH = $r }
    if ($targa -eq "GY983FY") { $rowGY983FY = $r }
}

# Record "before" operator values for the history log
$prevOperatorGL592TH = $wsStato.Cells.Item($rowGL592TH, 2).Value2
$prevOperatorGY983FY = $wsStato.Cells.Item($rowGY983FY, 2).Value2

# New operator / date values
$newOperatorGL592TH = "DI DEO MICHELE"
$newOperatorGY983FY = "FINE RENT"
$changeDate = "2026-02-13"

# Apply updates to "Stato Attuale"
$wsStato.Cells.Item($rowGL592TH, 2).Value = $newOperatorGL592TH
Set-TextValue $wsStato.Cells.Item($rowGL592TH, 3) $changeDate

$wsStato.Cells.Item($rowGY983FY, 2).Value = $newOperatorGY983FY
Set-TextValue $wsStato.Cells.Item($rowGY983FY, 3) $changeDate

# --- "Storico Passaggi" sheet: (re)write the two change-history rows ---
# The log only keeps the latest changes, so rows 2.. are overwritten fresh.
$row = 2

$wsStorico.Cells.Item($row, 1).Value = "GL592TH"
$wsStorico.Cells.Item($row, 2).Value = $prevOperatorGL592TH
$wsStorico.Cells.Item($row, 3).Value = $newOperatorGL592TH
Set-TextValue $wsStorico.Cells.Item($row, 4) $changeDate

$row = $row + 1

$wsStorico.Cells.Item($row, 1).Value = "GY983FY"
$wsStorico.Cells.Item($row, 2).Value = $prevOperatorGY983FY
$wsStorico.Cells.Item($row, 3).Value = $newOperatorGY983FY
Set-TextValue $wsStorico.Cells.Item($row, 4) $changeDate
